# Updated test cases sheet with shopping cart page (TC5 row)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 24) describing TC5 - "Check if user is able to proceed to checkout"
$ws.Range("A24").Value = "Shopping cart"
$ws.Range("B24").Value = "TC5"
$ws.Range("C24").Value = "Check if user is able to proceed to checkout"
$ws.Range("D24").Value = "High"
$ws.Range("E24").Value = "Regression"

# Move the active selection to E24 (matches the diff's sheetView selection change)
$ws.Range("E24").Select()
